# Grammar_REX.docx update (apollo 1/10/2018)
#  - scanner done / parser work in progress / token needs a better function
#
# This script rewrites a handful of grammar-table cells in the "Grammar"
# table using raw-OOXML paragraph replacement (InsertXML), which is the
# only way to get full control over run/bold boundaries and the
# w:proofErr / w:bookmarkStart|End markers that Find/Replace cannot touch.

$d = $word.ActiveDocument

function Replace-ParagraphXml($para, [string]$bodyXml) {
    # Replacing the *first* character position of an existing paragraph
    # with a full <w:p>...</w:p> fragment swaps out the whole paragraph
    # (pPr + all runs) in one shot - no leftover empty paragraphs.
    $insPoint = $para.Range
    $insPoint.Collapse(1) | Out-Null   # 1 = wdCollapseStart
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insPoint.InsertXML($pkg)
}

$tbl = $d.Tables(1)

# ---------------------------------------------------------------------
# Row "Function" (col 3): rule grew a parameter list and now ends in ';'
# ---------------------------------------------------------------------
$functionXml = '<w:p w:rsidR="002D2FB4" w:rsidRPr="00BF54CA" w:rsidRDefault="002D2FB4" w:rsidP="00AD7F94"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>FunctionName</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Type </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>VarName</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (,Type</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>VarName</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)*</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)?</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Block</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'

$functionCell = $tbl.Cell(4, 3)
Replace-ParagraphXml $functionCell.Range.Paragraphs(1) $functionXml

# ---------------------------------------------------------------------
# Row "InitializationTo0" (col 1): rule name gains "(oneDeclaration )"
# ---------------------------------------------------------------------
$initXml = '<w:p w:rsidR="002D2FB4" w:rsidRPr="00BF54CA" w:rsidRDefault="002D2FB4" w:rsidP="00AD7F94"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>InitializationTo0</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>oneDeclaration</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> )</w:t></w:r>' + `
  '</w:p>'

$d = $word.ActiveDocument
$initCell = $d.Tables(1).Cell(5, 1)
Replace-ParagraphXml $initCell.Range.Paragraphs(1) $initXml

# ---------------------------------------------------------------------
# Row "Comand" (col 3): paragraph mark becomes underlined, and the rule
# now allows repetition, so a trailing '*' run is appended.
# ---------------------------------------------------------------------
$comandXml = '<w:p w:rsidR="002D2FB4" w:rsidRDefault="00A41FB7" w:rsidP="00AD7F94"><w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Assignment</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>|LoopBlock|ConditionBlock|FunctionCall</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>|</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> end ;</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>|</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>giveBackWith</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Expression </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>|</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>giveBackWith</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> n</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>othing</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ;</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*</w:t></w:r>' + `
  '</w:p>'

$d = $word.ActiveDocument
$comandCell = $d.Tables(1).Cell(6, 3)
Replace-ParagraphXml $comandCell.Range.Paragraphs(1) $comandXml

# ---------------------------------------------------------------------
# Row "IfStatement" (col 3): drop the stray gramStart/gramEnd proofErr
# markers Word had put around "if(" - everything else is unchanged.
# ---------------------------------------------------------------------
$ifXml = '<w:p w:rsidR="00A41FB7" w:rsidRDefault="00A41FB7" w:rsidP="00A41FB7"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="009D2182"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>if(</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Operation</w:t></w:r>' + `
  '<w:r w:rsidRPr="009D2182"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>) start</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Comand</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>* (</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>el</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Operation</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> start </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Comand</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*)*(</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>else</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Comand</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*)?</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidRPr="009D2182"><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>fi</w:t></w:r>' + `
  '</w:p>'

$d = $word.ActiveDocument
$ifCell = $null
$tbl2 = $d.Tables(1)
for ($i = 1; $i -le $tbl2.Rows.Count; $i++) {
    $c = $tbl2.Cell($i, 3)
    if ($c.Range.Text -like "if(Operation*") {
        $ifCell = $c
        break
    }
}
Replace-ParagraphXml $ifCell.Range.Paragraphs(1) $ifXml

Write-Output "done"
